# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (GitHub Actions run on Fri May 26 13:23:42 UTC 2023).
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Rows 38/39 and 46/47 also swap which coin occupies which row (the
# source ranking re-ordered MXToken/VeChain and EnergySwap/Decentraland).
#
# D-column values that look numeric (e.g. "1.009", "99.00") are written
# with a leading apostrophe so Excel stores them as text instead of
# silently reformatting/rounding them - matching the workbook's existing
# inlineStr-based "Price" column, which intentionally keeps trailing
# zeros and dot-grouped big numbers (e.g. "26.645.02") as literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.645.02'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.830.86'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('D4').Value = '''1.009'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').Value = '''308.84'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = '''1.007'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('D7').Value = '''0.4681'
$ws.Range('E7').Value = '  +3.57%  '
$ws.Range('D8').Value = '''0.3605'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '''0.07146'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('D10').Value = '''0.9323'
$ws.Range('E10').Value = '  +4.66%  '
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('D12').Value = '''0.07646'
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('D13').Value = '1.821.86'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '''5.260'
$ws.Range('E14').Value = '  -0.78%  '
$ws.Range('D15').Value = '''6.353'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = '''87.65'
$ws.Range('E16').Value = '  +2.73%  '
$ws.Range('D17').Value = '''1.010'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').Value = '26.666.27'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').Value = '''14.27'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('D23').Value = '2.102.97'
$ws.Range('E23').Value = '  +2.53%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '''1.907'
$ws.Range('E25').Value = '  -3.10%  '
$ws.Range('D26').Value = '''151.70'
$ws.Range('E26').Value = '  +0.41%  '
$ws.Range('D28').Value = '''1.999'
$ws.Range('E28').Value = '  -2.70%  '
$ws.Range('D29').Value = '''113.75'
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('D30').Value = '''4.877'
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').Value = '''0.08829'
$ws.Range('E31').Value = '  +1.43%  '
$ws.Range('D32').Value = '''3.162'
$ws.Range('E32').Value = '  +1.42%  '
$ws.Range('D33').Value = '''2.845'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').Value = '''1.172'
$ws.Range('E34').Value = '  +5.24%  '
$ws.Range('D35').Value = '''0.7384'
$ws.Range('E35').Value = '  +1.86%  '
$ws.Range('D36').Value = '''4.446'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = '''1.079'
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.958'
$ws.Range('E38').Value = '  +2.24%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.01920'
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').Value = '''0.05150'
$ws.Range('E40').Value = '  +0.84%  '
$ws.Range('D41').Value = '''6.918'
$ws.Range('E41').Value = '  +1.68%  '
$ws.Range('D42').Value = '''0.5064'
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('D44').Value = '''8.113'
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''10.20'
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.4645'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('D48').Value = '''99.00'
$ws.Range('E48').Value = '  -1.86%  '
$ws.Range('D49').Value = '''1.574'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').Value = '''0.06027'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').Value = '''64.07'
$ws.Range('E51').Value = '  +0.17%  '
